$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# Date: 2025-05-21T14:22:51+00:00 -> 2025-06-13T15:45:04+00:00
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version: 4.3.0 -> 4.0.1
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet updates ---
# Row 2 (Extension) Constraint(s) - shorten the ele-1 constraint text
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 (Extension.id) Type(s): id -> string
$wsElem.Range("K3").Value = "string" + [char]10

# Row 6 (Extension.value[x]) Definition: R4B -> R4
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."

# Row 7 (Extension.value[x]:valueUrl) Definition: R4B -> R4
$wsElem.Range("M7").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
